$d = $word.ActiveDocument

# --- 1. Delete whole-paragraph ranges (do the later one first so the
#        earlier paragraph indices stay valid) ---

# Paragraphs 14-20: "If possible work out..." through the trailing
# empty paragraph after "To be decided!" (this also removes the
# bookmark that used to sit in "To be decided!").
$pStart = $d.Paragraphs(14)
$pEnd = $d.Paragraphs(20)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.Delete()

# Paragraphs 5-8: "Move crew positions..." through "Extend viewports..."
$pStart = $d.Paragraphs(5)
$pEnd = $d.Paragraphs(8)
$r = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$r.Delete()

# --- 2. Text replacements (order independent; Find scans the whole
#        story each time) ---

[void]$d.Content.Find.Execute("Morning", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do finding elevation of muzzle of gun", 2)

[void]$d.Content.Find.Execute("Restructure code I.E move all rendering operations into tank class", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do firing AP shell (including tracer model)", 2)

[void]$d.Content.Find.Execute("Make current operations work with new tank class", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do shell ballistics", 2)

[void]$d.Content.Find.Execute("Evening", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do changing seat function", 2)

[void]$d.Content.Find.Execute("Implement turning turret", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do proper UI", 2)

[void]$d.Content.Find.Execute("Remodel if necessary!", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do other vision blocks", 2)

[void]$d.Content.Find.Execute("Work out how to integrate collision objects in object file!", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Do tooltip for vision blocks", 2)

# --- 3. Insert a new empty paragraph right after "TODO!" ---
# (Using InsertAfter with a bare carriage return on a collapsed range
#  -- rather than InsertParagraphAfter -- keeps the new paragraph a
#  clean self-closing <w:p/> with no leftover empty run.)

$p1 = $d.Paragraphs(1)
$r1 = $p1.Range.Duplicate
$r1.Collapse(0)
$r1.InsertAfter([char]13)

# --- 4. Re-add the (now relocated) "_GoBack" bookmark as a zero-length
#        bookmark at the very start of the "Do finding elevation..."
#        paragraph (was paragraph 2, now paragraph 3 after the insert). ---

$pGun = $d.Paragraphs(3)
$bmRange = $d.Range($pGun.Range.Start, $pGun.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
